$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for the new (larger) account-statement table.
#    The old table spanned rows 16-41 (1 worker x up to 3 periods, 26 rows of
#    data). The new table spans rows 16-53 (38 rows of data). Insert 12 new
#    rows right after the old last data row (41) and before the footer rows
#    that used to be 46-47 (they will shift to 58-59).
# ---------------------------------------------------------------------------
$ws.Rows("42:53").Insert()

# Re-apply the correct cell formatting:
#  - row 53 becomes the new "last row" (it carries the bottom border), so
#    copy the formatting that used to belong to the old last row (41).
#  - rows 41-52 become ordinary data rows, so copy the formatting that
#    belongs to an ordinary data row (40).
$ws.Range("B41:J41").Copy()
$ws.Range("B53:J53").PasteSpecial(-4122)

$ws.Range("B40:J40").Copy()
$ws.Range("B41:J52").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Write out the new account-statement detail table (rows 16-53).
#    Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#             E=Periodo Mora, F=Valor Mora, G=Salario Basico
# ---------------------------------------------------------------------------
$rows = @(
  @(16,"CC","1098632764","EDINSON SAENZ ROJAS","2211",40000,1000000),
  @(17,"CC","1098632764","EDINSON SAENZ ROJAS","2209",40000,1000000),
  @(18,"CC","1098632764","EDINSON SAENZ ROJAS","2208",40000,1000000),
  @(19,"CC","3718491","OSCAR DE JESUS CAUSIL BURGOS","2211",40000,1000000),
  @(20,"CC","3718491","OSCAR DE JESUS CAUSIL BURGOS","2209",40000,1000000),
  @(21,"CC","3718491","OSCAR DE JESUS CAUSIL BURGOS","2208",40000,1000000),
  @(22,"CC","1050549637","OECSER DAMIAN MACHUCA DAMIAN","2211",40000,1000000),
  @(23,"CC","1050549637","OECSER DAMIAN MACHUCA DAMIAN","2209",40000,1000000),
  @(24,"CC","1050549637","OECSER DAMIAN MACHUCA DAMIAN","2208",40000,1000000),
  @(25,"CC","8865996","ADALBERTO ANTONIO CUELLO PEREZ","2211",40000,1000000),
  @(26,"CC","8865996","ADALBERTO ANTONIO CUELLO PEREZ","2209",40000,1000000),
  @(27,"CC","8865996","ADALBERTO ANTONIO CUELLO PEREZ","2208",40000,1000000),
  @(28,"CC","1093770388","HELVER ANDRES VARGAS DIAZ","2211",40000,1000000),
  @(29,"CC","1093770388","HELVER ANDRES VARGAS DIAZ","2209",40000,1000000),
  @(30,"CC","1093770388","HELVER ANDRES VARGAS DIAZ","2208",40000,1000000),
  @(31,"CC","91324739","BENJAMIN HEREDIA BARRAGAN","2211",40000,1000000),
  @(32,"CC","91324739","BENJAMIN HEREDIA BARRAGAN","2209",40000,1000000),
  @(33,"CC","91324739","BENJAMIN HEREDIA BARRAGAN","2208",40000,1000000),
  @(34,"CC","1102233009","LUZ ANGELA MIER RODRIGUEZ","2211",40000,1000000),
  @(35,"CC","1102233009","LUZ ANGELA MIER RODRIGUEZ","2209",40000,1000000),
  @(36,"CC","1102233009","LUZ ANGELA MIER RODRIGUEZ","2208",40000,1000000),
  @(37,"CC","1050548186","DIANA PATRICIA RICO SIERRA","2211",40000,1000000),
  @(38,"CC","1050548186","DIANA PATRICIA RICO SIERRA","2209",40000,1000000),
  @(39,"CC","1050548186","DIANA PATRICIA RICO SIERRA","2208",40000,1000000),
  @(40,"CC","1050555197","DEIVISON BALDOVINO ALEMAN","2211",40000,1000000),
  @(41,"CC","1050555197","DEIVISON BALDOVINO ALEMAN","2209",40000,1000000),
  @(42,"CC","1050555197","DEIVISON BALDOVINO ALEMAN","2208",40000,1000000),
  @(43,"CC","1104125815","MIYER SNEIDER CARPIO DIAZ","2211",40000,1000000),
  @(44,"CC","1104125815","MIYER SNEIDER CARPIO DIAZ","2209",40000,1000000),
  @(45,"CC","1104125815","MIYER SNEIDER CARPIO DIAZ","2208",40000,1000000),
  @(46,"CC","1002363066","LUIS EDUARDO GUERRERO GUERRERO","2303",37333,1000000),
  @(47,"CC","1002363066","LUIS EDUARDO GUERRERO GUERRERO","2302",40000,1000000),
  @(48,"CC","1002363066","LUIS EDUARDO GUERRERO GUERRERO","2301",40000,1000000),
  @(49,"CC","1002363066","LUIS EDUARDO GUERRERO GUERRERO","2212",40000,1000000),
  @(50,"CC","1002363066","LUIS EDUARDO GUERRERO GUERRERO","2211",40000,1000000),
  @(51,"CC","1002363066","LUIS EDUARDO GUERRERO GUERRERO","2210",40000,1000000),
  @(52,"CC","1002363066","LUIS EDUARDO GUERRERO GUERRERO","2209",40000,1000000),
  @(53,"CC","1002363066","LUIS EDUARDO GUERRERO GUERRERO","2208",40000,1000000)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 3. Update the summary figures above the table.
#    E11 = total "Valor Mora" (sum of column F for all workers/periods)
#    C13 = "Cant. Trabajadores" (distinct workers)
#    F13 = "Cant. Periodos" (max periods in arrears for a single worker)
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1517333
$ws.Range("C13").Value = 11
$ws.Range("F13").Value = 8
